# Iran weekly-deaths prediction workbook update.
# A new "day the prediction is made" block (2021-01-02) is inserted before
# the existing 2021-01-09 block, which gets pushed down by 10 rows
# (old rows 50-59 -> new rows 60-69, values unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Push the existing "2021-01-09" block (rows 50-59) down to rows 60-69.
$ws.Range("A50:K59").Insert(-4121)   # xlShiftDown

# Helper to write a literal text value into column A (dates like
# "2021-01-02" would otherwise be auto-recognised as real dates by Excel).
function Set-TextValue($cell, $text) {
    $cell.Value = "'" + $text
    $cell.ClearFormats()
}

# 2) New block: predictions made on 2021-01-02.
Set-TextValue $ws.Range("A50") "2021-01-02"
$ws.Range("B50").Value = "03 Jan -- 09 Jan 2021"
$ws.Range("C50").Value = 94.56999999999999
$ws.Range("D50").Value = 221.51
$ws.Range("E50").Value = 126.94
$ws.Range("F50").Value = "KNN"
$ws.Range("G50").Value = 2.89
$ws.Range("H50").Value = 48.92
$ws.Range("I50").Value = 60.38
$ws.Range("J50").Value = 83.98999999999999
$ws.Range("K50").Value = 85.28

Set-TextValue $ws.Range("A51") "2021-01-02"
$ws.Range("B51").Value = "10 Jan -- 16 Jan 2021"
$ws.Range("D51").Value = 245.27
$ws.Range("F51").Value = "KNN"

Set-TextValue $ws.Range("A52") "2021-01-02"
$ws.Range("B52").Value = "17 Jan -- 23 Jan 2021"
$ws.Range("D52").Value = 272.32
$ws.Range("F52").Value = "KNN"

Set-TextValue $ws.Range("A53") "2021-01-02"
$ws.Range("B53").Value = "24 Jan -- 30 Jan 2021"
$ws.Range("D53").Value = 300.93
$ws.Range("F53").Value = "KNN"

Set-TextValue $ws.Range("A54") "2021-01-02"
$ws.Range("B54").Value = "31 Jan -- 06 Feb 2021"
$ws.Range("D54").Value = 326.67
$ws.Range("F54").Value = "KNN"

Set-TextValue $ws.Range("A55") "2021-01-02"
$ws.Range("B55").Value = "07 Feb -- 13 Feb 2021"
$ws.Range("D55").Value = 320.26
$ws.Range("F55").Value = "KNN"

Set-TextValue $ws.Range("A56") "2021-01-02"
$ws.Range("B56").Value = "14 Feb -- 20 Feb 2021"
$ws.Range("D56").Value = 304.3
$ws.Range("F56").Value = "KNN"

Set-TextValue $ws.Range("A57") "2021-01-02"
$ws.Range("B57").Value = "21 Feb -- 27 Feb 2021"
$ws.Range("D57").Value = 291.37
$ws.Range("F57").Value = "KNN"

Set-TextValue $ws.Range("A58") "2021-01-02"
$ws.Range("B58").Value = "28 Feb -- 06 Mar 2021"
$ws.Range("D58").Value = 287.57
$ws.Range("F58").Value = "KNN"

Set-TextValue $ws.Range("A59") "2021-01-02"
$ws.Range("B59").Value = "07 Mar -- 13 Mar 2021"
$ws.Range("D59").Value = 273.42
$ws.Range("F59").Value = "KNN"
